# Update cryptos list cells based on the latest market data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.163.26"
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").Value = "2.245.16"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("D4").Value = "'1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.36%  "

$ws.Range("D5").Value = "'307.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.74%  "

$ws.Range("D6").Value = "'96.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.43%  "

$ws.Range("E7").Value = "  +1.27%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  -0.96%  "

$ws.Range("D10").Value = "'34.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.17%  "

$ws.Range("D11").Value = "'0.0817"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("D12").Value = "'7.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.96%  "

$ws.Range("E13").Value = "  +0.32%  "

$ws.Range("D14").Value = "2.587.43"
$ws.Range("E14").Value = "  +0.52%  "

$ws.Range("D15").Value = "2.244.95"
$ws.Range("E15").Value = "  +0.20%  "

$ws.Range("D16").Value = "'0.833"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("D17").Value = "'13.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.10%  "

$ws.Range("D18").Value = "44.064.72"
$ws.Range("E18").Value = "  +0.67%  "

$ws.Range("E19").Value = "  +1.59%  "

$ws.Range("D20").Value = "'12.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.32%  "

$ws.Range("D21").Value = "'6.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.95%  "

$ws.Range("D22").Value = "'65.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.14%  "

$ws.Range("D23").Value = "'236.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.77%  "

$ws.Range("E24").Value = "  -1.01%  "

$ws.Range("E25").Value = "  -1.30%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").Value = "'9.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.10%  "

$ws.Range("D28").Value = "'38.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.00%  "

$ws.Range("E29").Value = "  +1.54%  "

$ws.Range("D30").Value = "'6.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.27%  "

$ws.Range("D31").Value = "'20.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.96%  "

$ws.Range("D32").Value = "'152.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.34%  "

$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").Value = "'3.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.10%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0805"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.65%  "

$ws.Range("E35").Value = "  -3.21%  "

$ws.Range("E36").Value = "  +2.77%  "

$ws.Range("E37").Value = "  -1.03%  "

$ws.Range("D38").Value = "'1.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.15%  "

$ws.Range("D39").Value = "'14.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.59%  "

$ws.Range("D40").Value = "'3.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.16%  "

$ws.Range("E41").Value = "  -4.39%  "

$ws.Range("E42").Value = "  -2.61%  "

$ws.Range("D43").Value = "'1.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").Value = "1.723.66"
$ws.Range("E44").Value = "  +0.52%  "

$ws.Range("D45").Value = "'83.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.25%  "

$ws.Range("E46").Value = "  -1.00%  "

$ws.Range("D47").Value = "'100.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.96%  "

$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "'4.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.33%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'8.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.49%  "

$ws.Range("D50").Value = "'68.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.42%  "

$ws.Range("D51").Value = "'54.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.48%  "
